# Apply updated crypto price/volume data to sheet1 (matches commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.200.03'
$ws.Range('D3').Value = '2.302.99'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'516.89"
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = "'130.84"
$ws.Range('E6').Value = '  -3.65%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D9').Value = '2.316.66'
$ws.Range('E9').Value = '  -1.22%  '
$ws.Range('D10').Value = "'0.0995"
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').Value = "'0.153"
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = "'5.24"
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('D13').Value = "'0.336"
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '2.716.91'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = "'23.21"
$ws.Range('E15').Value = '  -3.45%  '
$ws.Range('D16').Value = '56.148.69'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('E17').Value = '  -2.62%  '
$ws.Range('D18').Value = '2.293.14'
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('D19').Value = "'329.58"
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').Value = "'10.31"
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('D21').Value = "'4.12"
$ws.Range('E21').Value = '  -2.75%  '
$ws.Range('D22').Value = "'6.68"
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = "'60.88"
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = "'8.58"
$ws.Range('E25').Value = '  +7.15%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = "'0.163"
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('D27').Value = "'0.994"
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = "'1.31"
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D29').Value = "'167.95"
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('E31').Value = '  -4.63%  '
$ws.Range('D32').Value = "'6.07"
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('D33').Value = "'18.20"
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D35').Value = "'0.996"
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('E36').Value = '  -2.58%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = "'3.88"
$ws.Range('E37').Value = '  -3.37%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').Value = "'0.880"
$ws.Range('E38').Value = '  -4.39%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').Value = "'38.54"
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('D41').Value = "'147.78"
$ws.Range('E41').Value = '  +3.81%  '
$ws.Range('D42').Value = "'0.372"
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('D43').Value = "'283.69"
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').Value = "'5.04"
$ws.Range('E45').Value = '  -4.08%  '
$ws.Range('D46').Value = "'0.0924"
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('D47').Value = "'0.0495"
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').Value = "'0.554"
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('D49').Value = "'18.12"
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('E51').Value = '  -2.97%  '
